$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Update M3 and N3 values
$ws.Range("M3").Value = "عطل"
$ws.Range("N3").Value = "تم اصلاح "

# Delete column O entirely (removes header "Serviced by " and empty cells below)
$ws.Range("O1:O12").EntireColumn.Delete()
